$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-08 Saturday" "2025-11-09 Sunday"

Replace-Text "76×77=5852" "77×51=3927"
Replace-Text "50×55=2750" "87×44=3828"
Replace-Text "26×55=1430" "25×30=750"
Replace-Text "72×97=6984" "70×32=2240"
Replace-Text "91×38=3458" "67×31=2077"

Replace-Text "51×13=663" "39×60=2340"
Replace-Text "46×86=3956" "47×29=1363"
Replace-Text "27×26=702" "81×30=2430"
Replace-Text "58×27=1566" "16×57=912"
Replace-Text "16×32=512" "28×97=2716"

Replace-Text "30×23=690" "46×68=3128"
Replace-Text "94×67=6298" "92×67=6164"
Replace-Text "56×31=1736" "13×45=585"
Replace-Text "15×38=570" "85×41=3485"
Replace-Text "85×70=5950" "95×90=8550"

Replace-Text "97×19=1843" "87×47=4089"
Replace-Text "21×17=357" "27×35=945"
Replace-Text "21×85=1785" "92×56=5152"
Replace-Text "96×42=4032" "88×86=7568"
Replace-Text "20×50=1000" "83×33=2739"

Replace-Text "85×75=6375" "50×45=2250"
Replace-Text "85×81=6885" "56×61=3416"
Replace-Text "27×77=2079" "74×71=5254"
Replace-Text "74×12=888" "63×53=3339"
Replace-Text "54×71=3834" "83×16=1328"
